$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.295.76"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.708.66"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5298"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2641"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06547"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.90"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07645"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.704.26"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.945.13"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5744"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8177"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.28"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.283.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.86"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.669"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.962"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.43"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.755"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1216"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.268"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05369"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.294"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.483"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.415"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.879"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9536"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.422"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5863"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01625"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.880"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.042.56"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8393"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.02"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.851.76"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4500"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.053"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06506"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.50%  "
